$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NextBus1")

# Row 2
$ws.Range("F2").Value = 45702.46756944444
$ws.Range("L2").Value = "DD"
$ws.Range("O2").Value = 16

# Row 3
$ws.Range("F3").Value = 45702.45746527778
$ws.Range("O3").Value = 1

# Row 4
$ws.Range("F4").Value = 45702.45915509259
$ws.Range("O4").Value = 4

# Row 5
$ws.Range("F5").Value = 45702.46603009259
$ws.Range("O5").Value = 14

# Row 6
$ws.Range("F6").Value = 45702.45571759259
$ws.Range("O6").Value = 0

# Row 7
$ws.Range("F7").Value = 45702.46324074074
$ws.Range("L7").Value = "SD"
$ws.Range("O7").Value = 10

# Row 8
$ws.Range("F8").Value = 45702.46046296296
$ws.Range("I8").Value = "SEA"
$ws.Range("O8").Value = 6

# Row 9
$ws.Range("F9").Value = 45702.45592592593

# Row 10
$ws.Range("F10").Value = 45702.45571759259
$ws.Range("O10").Value = 0

# Row 11
$ws.Range("F11").Value = 45702.45809027777
$ws.Range("O11").Value = 2

# Row 12
$ws.Range("F12").Value = 45702.45716435185
$ws.Range("L12").Value = "SD"
$ws.Range("O12").Value = 1

# Row 13
$ws.Range("F13").Value = 45702.46082175926
$ws.Range("O13").Value = 6

# Row 14
$ws.Range("F14").Value = 45702.46027777778
$ws.Range("O14").Value = 6

# Row 15
$ws.Range("F15").Value = 45702.46402777778
$ws.Range("O15").Value = 11
